# Update the "L1" sheet (Buffer Used bytes / Packets Dropped columns)
# and the "F1" sheet (Window Size bps column) with newly-simulated
# Reno congestion-control numbers.

$wb = $excel.ActiveWorkbook

$wsL1 = $wb.Worksheets.Item("L1")
$wsF1 = $wb.Worksheets.Item("F1")

# --- L1 sheet: column C (Buffer Used), rows 4-51 get new values; row 3 unchanged ---
$l1C = @{
    4  = 4480
    5  = 11520
    6  = 9600
    7  = 19200
    8  = 16000
    9  = 28800
    10 = 22400
    11 = 38400
    12 = 28160
    13 = 46080
    14 = 33920
    15 = 55680
    16 = 40320
    17 = 65280
    18 = 46720
    19 = 74880
    20 = 53120
    21 = 84480
    22 = 59520
    23 = 94080
    24 = 65920
    25 = 103680
    26 = 72320
    27 = 113280
    28 = 78720
    29 = 39680
    30 = 0
    31 = 0
    32 = 0
    33 = 0
    34 = 0
    35 = 0
    36 = 0
    37 = 0
    38 = 0
    39 = 0
    40 = 0
    41 = 0
    42 = 0
    43 = 0
    44 = 0
    45 = 0
    46 = 0
    47 = 0
    48 = 0
    49 = 0
    50 = 0
    51 = 0
}

foreach ($row in $l1C.Keys) {
    $wsL1.Cells.Item($row, 3).Value = $l1C[$row]
}

# --- L1 sheet: column E (Packets Dropped), rows 3-51 all become 0 ---
for ($row = 3; $row -le 51; $row++) {
    $wsL1.Cells.Item($row, 5).Value = 0
}

# --- F1 sheet: column C (Window Size bps), rows 2-51 get new values ---
$f1C = @{
    2  = 1
    3  = 2.9
    4  = 5.197631445038131
    5  = 7.02286464978477
    6  = 9.676964161540262
    7  = 11.48108737366293
    8  = 14.2840640402846
    9  = 16.19135853865013
    10 = 19.08729694930582
    11 = 21.03383311293605
    12 = 23.92804176884484
    13 = 25.81912680736002
    14 = 28.6833153657839
    15 = 30.60714843611165
    16 = 33.50952365744345
    17 = 35.4534143544298
    18 = 38.37999394672695
    19 = 40.33700142911891
    20 = 43.27983505182375
    21 = 45.24586785540937
    22 = 48.20012264220149
    23 = 50.17262140792511
    24 = 53.13519850160618
    25 = 55.11248372847466
    26 = 58.08130799627527
    27 = 60.06223339773678
    28 = 63.03586447177104
    29 = 64.0119955810089
    30 = 64.0119955810089
    31 = 64.0119955810089
    32 = 64.0119955810089
    33 = 64.0119955810089
    34 = 64.0119955810089
    35 = 64.0119955810089
    36 = 64.0119955810089
    37 = 64.0119955810089
    38 = 64.0119955810089
    39 = 64.0119955810089
    40 = 64.0119955810089
    41 = 64.0119955810089
    42 = 64.0119955810089
    43 = 64.0119955810089
    44 = 64.0119955810089
    45 = 64.0119955810089
    46 = 64.0119955810089
    47 = 64.0119955810089
    48 = 64.0119955810089
    49 = 64.0119955810089
    50 = 64.0119955810089
    51 = 64.0119955810089
}

foreach ($row in $f1C.Keys) {
    $wsF1.Cells.Item($row, 3).Value = $f1C[$row]
}
